$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the trailing three data rows (old "MuSCs" sending-cluster block, rows 8-10)
$ws.Range("A8:A10").EntireRow.Delete()

# Update remaining data rows (2-7) with the new TPM-derived values
# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Fgf7"
$ws.Range("C2").Value = "Fgfr3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 17.39906333333333
$ws.Range("H2").Value = 52.19719000000001
$ws.Range("I2").Value = 0.9351306508759385
$ws.Range("J2").Value = 0.9351306508759385
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 4.959409333333333
$ws.Range("N2").Value = 14.878228
$ws.Range("O2").Value = 0.8271666313262851
$ws.Range("P2").Value = 0.8271666313262852
$ws.Range("Q2").Value = 86.28907708659112
$ws.Range("R2").Value = 776.60169377932
$ws.Range("S2").Value = 0.7735088703350065
$ws.Range("T2").Value = 0.7735088703350066

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Fgf7"
$ws.Range("C3").Value = "Fgfr3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 17.39906333333333
$ws.Range("H3").Value = 52.19719000000001
$ws.Range("I3").Value = 0.9351306508759385
$ws.Range("J3").Value = 0.9351306508759385
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.5648773333333333
$ws.Range("N3").Value = 1.694632
$ws.Range("O3").Value = 0.09421438109281059
$ws.Range("P3").Value = 0.09421438109281059
$ws.Range("Q3").Value = 9.828336498231112
$ws.Range("R3").Value = 88.45502848408
$ws.Range("S3").Value = 0.08810275551319369
$ws.Range("T3").Value = 0.08810275551319369

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Fgf7"
$ws.Range("C4").Value = "Fgfr3"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 17.39906333333333
$ws.Range("H4").Value = 52.19719000000001
$ws.Range("I4").Value = 0.9351306508759385
$ws.Range("J4").Value = 0.9351306508759385
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.4713726666666667
$ws.Range("N4").Value = 1.414118
$ws.Range("O4").Value = 0.07861898758090437
$ws.Range("P4").Value = 0.07861898758090438
$ws.Range("Q4").Value = 8.201442880935556
$ws.Range("R4").Value = 73.81298592842
$ws.Range("S4").Value = 0.07351902502773842
$ws.Range("T4").Value = 0.07351902502773844

# Row 5
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Fgf7"
$ws.Range("C5").Value = "Fgfr3"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.206960666666667
$ws.Range("H5").Value = 3.620882
$ws.Range("I5").Value = 0.06486934912406146
$ws.Range("J5").Value = 0.06486934912406146
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 4.959409333333333
$ws.Range("N5").Value = 14.878228
$ws.Range("O5").Value = 0.8271666313262851
$ws.Range("P5").Value = 0.8271666313262852
$ws.Range("Q5").Value = 5.985811995232888
$ws.Range("R5").Value = 53.872307957096
$ws.Range("S5").Value = 0.05365776099127862
$ws.Range("T5").Value = 0.05365776099127863

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Fgf7"
$ws.Range("C6").Value = "Fgfr3"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.206960666666667
$ws.Range("H6").Value = 3.620882
$ws.Range("I6").Value = 0.06486934912406146
$ws.Range("J6").Value = 0.06486934912406146
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.5648773333333333
$ws.Range("N6").Value = 1.694632
$ws.Range("O6").Value = 0.09421438109281059
$ws.Range("P6").Value = 0.09421438109281059
$ws.Range("Q6").Value = 0.6817847228248889
$ws.Range("R6").Value = 6.136062505423999
$ws.Range("S6").Value = 0.006111625579616906
$ws.Range("T6").Value = 0.006111625579616906

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Fgf7"
$ws.Range("C7").Value = "Fgfr3"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.206960666666667
$ws.Range("H7").Value = 3.620882
$ws.Range("I7").Value = 0.06486934912406146
$ws.Range("J7").Value = 0.06486934912406146
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.4713726666666667
$ws.Range("N7").Value = 1.414118
$ws.Range("O7").Value = 0.07861898758090437
$ws.Range("P7").Value = 0.07861898758090438
$ws.Range("Q7").Value = 0.5689282680084444
$ws.Range("R7").Value = 5.120354412076
$ws.Range("S7").Value = 0.005099962553165938
$ws.Range("T7").Value = 0.005099962553165939
